# Generate Report for Archive
#
# The localization status for two files (6f37c1fc-...md and
# e7c82c9c-...md) moved from "Ready for handoff" back to "In Translation"
# on both locales (zh-cn, de-de); the Overview sheet mirrors the same
# Status values for each locale column. The third file in that block
# (fd38aa26-...md) keeps its "Ready for handoff" status.

$wb = $excel.ActiveWorkbook

$inTranslation = "In Translation"

# --- Overview sheet: columns B (zh-cn) and C (de-de) mirror status ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B8").Value = $inTranslation
$overview.Range("C8").Value = $inTranslation
$overview.Range("B9").Value = $inTranslation
$overview.Range("C9").Value = $inTranslation

# --- zh-cn sheet: column C is Status ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C8").Value = $inTranslation
$zhcn.Range("C9").Value = $inTranslation

# --- de-de sheet: column C is Status ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C8").Value = $inTranslation
$dede.Range("C9").Value = $inTranslation
